# Q&As.xlsx - add two new Q/A/Reasoning rows (rows 5 and 6) with sentence
# embedding context pasted from AR50 / matrikkel documentation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: "Videregående skole" distance question -------------------------
# Write order matches how the shared strings ended up appended in the
# original commit: Question (A5), Answer (B5), Reasoning (C5).
$ws.Range("A5").Value2 = "How many percent of buildings have more than 2km to the nearest ""Videregående skole""?"
$ws.Range("B5").Value2 = "~31%"
$ws.Range("C5").Value2 = "Steps:
1) Gather and load building point  data and its documentation
2) Find the ""bygningstype"" field and check the documentation to see that ""Videregående skole"" has type 616
3) Select building points that have ""bygninigstype"" = 616
4) Create a 2-kilometer buffer around the selected points
5) Optionally dissolve the buffers
6) Perform an difference or intersection operation to separate points that are inside/outside the buffer(s)
7) Perform the percentage calculations"

# --- Row 6: forest productivity question ------------------------------------
# Write order: Answer (B6) first, then Question (A6), then Reasoning (C6),
# matching the shared-string insertion order seen in the target workbook.
$ws.Range("B6").Value2 = "53,855283 km^2"
$ws.Range("A6").Value2 = "Determine how many square kilometers have high forest productivity."
$ws.Range("C6").Value2 = "Steps:
1) Gather and load building point AR50 data , either the entire series or only ""Jordbruk"", along with the AR50 documentation
2) Find the ""skogbonitet"" field and check the documentation to see that areas with the highest forest productivity have value 18.
3) Select building points that have ""skogbonitet"" = 18
4) Calculate the area of each selected polygon.
5) Sum all areas."

# B6 was pasted in from elsewhere (e.g. a web page), carrying its own font
# (same Calibri 11, but an explicit dark-gray color) instead of the sheet's
# default font/theme color, while keeping the column's left/right border.
$srcFormat = $ws.Range("B5")
$srcFormat.Copy()
$ws.Range("B6").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("B6").WrapText = $false
$ws.Range("B6").Font.Color = 2367776

# --- Row heights (content no longer fits the default 14.4pt row) -----------
$ws.Rows.Item(5).RowHeight = 172.8
$ws.Rows.Item(6).RowHeight = 144

# --- View state --------------------------------------------------------------
$ws.Range("B9").Select()

Write-Host "Applied Q&A rows 5-6 edits."
